$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the crypto price refresh.
# D-column cells receive special text-forcing handling because some
# values look numeric (e.g. "1.443", "0.9993") and Excel would
# otherwise auto-convert them to numbers instead of keeping them as text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.440.26"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.666.99"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.17"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("E7").Value = "  +1.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3913"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.08"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").Value = "  +5.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.405"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = "  +2.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9992"
$ws.Range("D11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08582"
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.48"
$ws.Range("D13").ClearFormats()

$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.325"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = "  +1.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001341"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").Value = "  +3.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.887"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").Value = "  +4.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.662.72"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = "  +0.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.52"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06965"
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.54"
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = "  -3.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.995"
$ws.Range("D21").ClearFormats()

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9977"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.71"
$ws.Range("D23").ClearFormats()

$ws.Range("E23").Value = "  -1.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.428.79"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = "  +0.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.434"
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = "  +2.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.028"
$ws.Range("D26").ClearFormats()

$ws.Range("E26").Value = "  +9.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.51"
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.98"
$ws.Range("D28").ClearFormats()

$ws.Range("E28").Value = "  -0.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "142.77"
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.347"
$ws.Range("D30").ClearFormats()

$ws.Range("E30").Value = "  -0.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.096"
$ws.Range("D31").ClearFormats()

$ws.Range("E31").Value = "  -6.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.504"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = "  +4.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.843.53"
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.068"
$ws.Range("D34").ClearFormats()

$ws.Range("E34").Value = "  +6.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08267"
$ws.Range("D35").ClearFormats()

$ws.Range("E35").Value = "  +1.44%  "

$ws.Range("B36").Value = "VeChain"

$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02985"
$ws.Range("D36").ClearFormats()

$ws.Range("E36").Value = "  +1.39%  "

$ws.Range("B37").Value = "FraxShare"

$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.23"
$ws.Range("D37").ClearFormats()

$ws.Range("E37").Value = "  +10.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.804"
$ws.Range("D38").ClearFormats()

$ws.Range("E38").Value = "  -3.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2757"
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = "  +1.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09260"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7737"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.82"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = "  +4.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.445"
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = "  -2.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.55"
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = "  +2.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7111"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = "  +2.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.525"
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.141"
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = "  +0.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9975"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = "  -0.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08465"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.33"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").Value = "  +0.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.443"
$ws.Range("D51").ClearFormats()

$ws.Range("E51").Value = "  +11.80%  "
